$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 01:51"

# Update country data rows: new scraped totals shuffled the ranking
# (table is sorted by total cases, column B, descending), so some
# rows both get new numbers AND a different country name.

# Row 4: 'Estados Unidos' (data refresh)
$ws.Range("B4").Value = 7545793
$ws.Range("C4").Value = 47873
$ws.Range("D4").Value = 4770317
$ws.Range("E4").Value = 2562025
$ws.Range("G4").Value = 791
$ws.Range("H4").Value = 213451

# Row 6: 'Brasil' (data refresh)
$ws.Range("B6").Value = 4882231
$ws.Range("C6").Value = 33002
$ws.Range("D6").Value = 4232593
$ws.Range("E6").Value = 504207
$ws.Range("G6").Value = 664
$ws.Range("H6").Value = 145431

# Row 29: 'Canada' (data refresh)
$ws.Range("B29").Value = 162659
$ws.Range("C29").Value = 2124
$ws.Range("D29").Value = 137614
$ws.Range("E29").Value = 15636
$ws.Range("G29").Value = 90
$ws.Range("H29").Value = 9409

# Row 37: 'Republica Dominicana' -> 'Panama'
$ws.Range("A37").Value = "Panama"
$ws.Range("B37").Value = 113962
$ws.Range("C37").Value = 620
$ws.Range("D37").Value = 90772
$ws.Range("E37").Value = 20784
$ws.Range("G37").Value = 19
$ws.Range("H37").Value = 2406

# Row 38: 'Panama' -> 'Republica Dominicana'
$ws.Range("A38").Value = "Republica Dominicana"
$ws.Range("B38").Value = 113350
$ws.Range("C38").Value = 622
$ws.Range("D38").Value = 88840
$ws.Range("E38").Value = 22393
$ws.Range("G38").Value = 9
$ws.Range("H38").Value = 2117

# Row 51: 'Costa Rica' -> 'Chequia'
$ws.Range("A51").Value = "Chequia"
$ws.Range("B51").Value = 78051
$ws.Range("C51").Value = 3796
$ws.Range("D51").Value = 35032
$ws.Range("E51").Value = 42320
$ws.Range("G51").Value = 21
$ws.Range("H51").Value = 699

# Row 52: 'Honduras' -> 'Costa Rica'
$ws.Range("A52").Value = "Costa Rica"
$ws.Range("B52").Value = 77829
$ws.Range("C52").Value = 1001
$ws.Range("D52").Value = 42621
$ws.Range("E52").Value = 34278
$ws.Range("G52").Value = 13
$ws.Range("H52").Value = 930

# Row 53: 'Portugal' -> 'Honduras'
$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 77598
$ws.Range("C53").Value = 698
$ws.Range("D53").Value = 28517
$ws.Range("E53").Value = 46701
$ws.Range("G53").Value = 27
$ws.Range("H53").Value = 2380

# Row 54: 'Etiopia' -> 'Portugal'
$ws.Range("A54").Value = "Portugal"
$ws.Range("B54").Value = 77284
$ws.Range("C54").Value = 888
$ws.Range("D54").Value = 49359
$ws.Range("E54").Value = 25942
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 1983

# Row 55: 'Venezuela' -> 'Etiopia'
$ws.Range("A55").Value = "Etiopia"
$ws.Range("B55").Value = 76988
$ws.Range("C55").Value = 890
$ws.Range("D55").Value = 31677
$ws.Range("E55").Value = 44103
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 1208

# Row 56: 'Chequia' -> 'Venezuela'
$ws.Range("A56").Value = "Venezuela"
$ws.Range("B56").Value = 76029
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 66245
$ws.Range("E56").Value = 9149
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 635

# Row 81: 'Australia' (data refresh)
$ws.Range("B81").Value = 27113
$ws.Range("C81").Value = 17
$ws.Range("D81").Value = 24824
$ws.Range("E81").Value = 1399

# Row 84: 'Camerun' (data refresh)
$ws.Range("B84").Value = 20924
$ws.Range("D84").Value = 19764
$ws.Range("E84").Value = 740
$ws.Range("H84").Value = 420

# Row 94: 'Noruega' (data refresh)
$ws.Range("B94").Value = 14284
$ws.Range("C94").Value = 135
$ws.Range("E94").Value = 2819

# Row 102: 'Consejo Danes para los Refugiados' -> 'Guinea'
$ws.Range("A102").Value = "Guinea"
$ws.Range("B102").Value = 10735
$ws.Range("C102").Value = 83
$ws.Range("D102").Value = 10066
$ws.Range("E102").Value = 603
$ws.Range("H102").Value = 66

# Row 103: 'Guinea' -> 'Consejo Danes para los Refugiados'
$ws.Range("A103").Value = "Consejo Danes para los Refugiados"
$ws.Range("B103").Value = 10729
$ws.Range("C103").Value = 44
$ws.Range("D103").Value = 10183
$ws.Range("E103").Value = 274
$ws.Range("H103").Value = 272

# Row 130: 'Surinam' (data refresh)
$ws.Range("B130").Value = 4899
$ws.Range("C130").Value = 8
$ws.Range("D130").Value = 4715
$ws.Range("E130").Value = 79

# Row 136: 'Aruba' -> 'Reunion'
$ws.Range("A136").Value = "Reunion"
$ws.Range("B136").Value = 4178
$ws.Range("C136").Value = 185
$ws.Range("D136").Value = 3360
$ws.Range("E136").Value = 802
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 16

# Row 137: 'Reunion' -> 'Aruba'
$ws.Range("A137").Value = "Aruba"
$ws.Range("B137").Value = 4038
$ws.Range("C137").Value = 40
$ws.Range("D137").Value = 3406
$ws.Range("E137").Value = 602
$ws.Range("G137").Value = 3
$ws.Range("H137").Value = 30

# Row 155: 'Uruguay' (data refresh)
$ws.Range("B155").Value = 2097
$ws.Range("C155").Value = 36
$ws.Range("D155").Value = 1824
$ws.Range("E155").Value = 225

# Row 183: 'Eritrea' -> 'Mauricio'
$ws.Range("A183").Value = "Mauricio"
$ws.Range("B183").Value = 385
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 344
$ws.Range("E183").Value = 31
$ws.Range("H183").Value = 10

# Row 184: 'Mauricio' -> 'Eritrea'
$ws.Range("A184").Value = "Eritrea"
$ws.Range("D184").Value = 353
$ws.Range("E184").Value = 28
$ws.Range("H184").Value = 0

# Row 190: 'Islas Caimanes' (data refresh)
$ws.Range("B190").Value = 213
$ws.Range("C190").Value = 2
$ws.Range("E190").Value = 2
